# Updated cryptos list on Tue Jun  6 16:41:17 UTC 2023 with GitHub Actions
# Refresh per-coin Price (D) / Volume(1h) (E) values; two coin rows swapped
# (VeChain/TheSandbox at 39-40, Cronos/NEARProtocol at 50-51).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric-looking Price/Volume text must stay text (apostrophe-prefix
# forces literal-text entry like a user typing it in Excel), then reset
# the cell style back to Normal so no stray quote-prefix/number-format
# style sticks to the cell.

$ws.Range("D2").Value = "'26.163.43"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.92%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.843.13"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +1.22%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'0.9998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.51%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'279.35"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.91%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'0.9994"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -0.45%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.5093"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +0.71%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.3505"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -1.11%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'45.00"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -0.12%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.06830"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +2.12%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'19.99"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -0.79%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.8056"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -5.38%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.07782"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -1.20%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'1.849.14"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +1.50%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'5.093"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +0.98%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'88.53"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +1.31%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.9995"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -0.62%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'14.19"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.94%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.000008061"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -0.88%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.9990"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -0.43%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'26.211.30"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +0.82%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'4.771"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.38%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'10.08"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -0.37%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'6.211"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +1.41%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'2.383"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +10.69%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'144.52"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +2.05%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'1.666"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -0.87%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'17.22"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +1.40%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'110.07"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +1.18%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'4.367"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +1.39%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'4.299"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +1.39%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'0.08742"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -0.78%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.04908"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +2.41%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'1.169"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +3.82%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'0.7341"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -1.08%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'2.840"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -0.35%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'3.232"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +4.25%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'2.386"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -1.68%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").Value = "'0.01849"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -0.46%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("B40").Value = 'TheSandbox'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D40").Value = "'0.5171"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -4.42%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.9656"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -2.08%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'116.24"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +3.40%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'6.269"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +0.70%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'8.002"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'0.9985"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -0.55%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.4524"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -4.26%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.1356"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -1.50%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'9.326"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +0.66%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  +1.62%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").Value = "'0.05932"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +0.24%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("B51").Value = 'NEARProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D51").Value = "'1.501"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +0.27%  "
$ws.Range("E51").Style = "Normal"
